# Q3 2024 Fiscal Update
# Adds Jul-Oct 2024 (columns L:O) monthly BTr NG Cash Operations data to the
# "2024" sheet, corrects a handful of already-reported Jan-Jun 2024 figures,
# and un-hides the columns that previously held no data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# ---------------------------------------------------------------------
# 1) Un-hide columns J:Q (J:K were already visible/used; L:Q held no data
#    and were hidden -- now Jul-Oct data populates L:O so reveal them all)
# ---------------------------------------------------------------------
$ws.Range("J1:Q1").EntireColumn.Hidden = $false

# ---------------------------------------------------------------------
# 2) Corrections to existing Jan-Jun 2024 (cols F:K) figures
# ---------------------------------------------------------------------
$ws.Range("F19").Value = 3328
$ws.Range("G19").Value = 2762
$ws.Range("H19").Value = 3635
$ws.Range("I19").Value = 2507
$ws.Range("J19").Value = 3366
$ws.Range("K19").Value = 2186

$ws.Range("F23").Value = 4570
$ws.Range("G23").Value = 3665
$ws.Range("H23").Value = 3871
$ws.Range("I23").Value = 4852
$ws.Range("J23").Value = 3517
$ws.Range("K23").Value = 1417

$ws.Range("F26").Value = 14156
$ws.Range("G26").Value = 1825
$ws.Range("H26").Value = 10588
$ws.Range("I26").Value = 6113
$ws.Range("J26").Value = 53181
$ws.Range("K26").Value = 36814

# Rows that used to compute via formula from other cells on the sheet now
# become plain reported figures (the figures themselves are unchanged).
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 1587
$ws.Range("H33").Value = 4811
$ws.Range("I33").Value = 6911
$ws.Range("J33").Value = 1015

$ws.Range("F37").Value = 181265
$ws.Range("G37").Value = 225033
$ws.Range("H37").Value = 321416
$ws.Range("I37").Value = 311729
$ws.Range("J37").Value = 403765
$ws.Range("K37").Value = 349713

$ws.Range("F39").Value = 87951
$ws.Range("G39").Value = -164677
$ws.Range("H39").Value = -195918
$ws.Range("I39").Value = 42728
$ws.Range("J39").Value = -174911
$ws.Range("K39").Value = -209080

$ws.Range("F42").Value = -22893
$ws.Range("G42").Value = 2578
$ws.Range("H42").Value = 44201
$ws.Range("I42").Value = -32260
$ws.Range("J42").Value = 119815
$ws.Range("K42").Value = 7845

# ---------------------------------------------------------------------
# 3) New Jul-Oct 2024 (cols L:O) data
# ---------------------------------------------------------------------

# Row 11 - BIR
$ws.Range("L11").Value = 319814
$ws.Range("M11").Value = 238120
$ws.Range("N11").Value = 174679
$ws.Range("O11").Value = 325536

# Row 13 - Documentary Stamp
$ws.Range("L13").Value = 1258
$ws.Range("M13").Value = 1272
$ws.Range("N13").Value = 3326
$ws.Range("O13").Value = 0

# Row 14 - Tax Expenditures (BIR)
$ws.Range("L14").Value = 380
$ws.Range("M14").Value = 12
$ws.Range("N14").Value = 2828
$ws.Range("O14").Value = 362

# Row 16 - BOC
$ws.Range("L16").Value = 80355
$ws.Range("M16").Value = 78521
$ws.Range("N16").Value = 76282
$ws.Range("O16").Value = 86885

# Row 18 - Tax Expenditures (BOC)
$ws.Range("L18").Value = 2675
$ws.Range("M18").Value = 814
$ws.Range("N18").Value = 1001
$ws.Range("O18").Value = 566

# Row 19 - Other Offices
$ws.Range("L19").Value = 2653
$ws.Range("M19").Value = 3572
$ws.Range("N19").Value = 2498
$ws.Range("O19").Value = 2446

# Row 22 - BTr Income
$ws.Range("L22").Value = 19914
$ws.Range("M22").Value = 16485
$ws.Range("N22").Value = 9921
$ws.Range("O22").Value = 14519

# Row 23 - Fees and Charges
$ws.Range("L23").Value = 1397
$ws.Range("M23").Value = 1487
$ws.Range("N23").Value = 1421
$ws.Range("O23").Value = 877

# Row 24 - Privatization
$ws.Range("L24").Value = 18
$ws.Range("M24").Value = 2684
$ws.Range("N24").Value = 6
$ws.Range("O24").Value = 146

# Row 25 - Income from Malampaya
$ws.Range("L25").Value = 965
$ws.Range("M25").Value = 1104
$ws.Range("N25").Value = 785
$ws.Range("O25").Value = 1168

# Row 26 - Other non-tax
$ws.Range("L26").Value = 32257
$ws.Range("M26").Value = 44204
$ws.Range("N26").Value = 34058
$ws.Range("O26").Value = 41547

# Row 28 - Grants
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 160
$ws.Range("N28").Value = 3
$ws.Range("O28").Value = 0

# Row 31 - Allotment to LGUs
$ws.Range("L31").Value = 86190
$ws.Range("M31").Value = 82515
$ws.Range("N31").Value = 82505
$ws.Range("O31").Value = 84378

# Row 32 - Interest Payments
$ws.Range("L32").Value = 79432
$ws.Range("M32").Value = 52781
$ws.Range("N32").Value = 73852
$ws.Range("O32").Value = 55388

# Row 33 - Tax Expenditures
$ws.Range("L33").Value = 4313
$ws.Range("M33").Value = 2098
$ws.Range("N33").Value = 7155
$ws.Range("O33").Value = 928

# Row 34 - Equity
$ws.Range("L34").Value = 10719
$ws.Range("M34").Value = 9100
$ws.Range("N34").Value = 18217
$ws.Range("O34").Value = 11968

# Row 35 - Net Lending
$ws.Range("L35").Value = 85
$ws.Range("M35").Value = 1
$ws.Range("N35").Value = 381
$ws.Range("O35").Value = 22

# Row 36 - NG Disbursements
$ws.Range("L36").Value = -1143
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 3115
$ws.Range("O36").Value = 3053

# Row 37 - Subsidy
$ws.Range("L37").Value = 306623
$ws.Range("M37").Value = 294048
$ws.Range("N37").Value = 387694
$ws.Range("O37").Value = 311047

# Row 39 - Surplus/(-)Deficit
$ws.Range("L39").Value = -28845
$ws.Range("M39").Value = -54206
$ws.Range("N39").Value = -273266
$ws.Range("O39").Value = 6340

# Row 42 - External (Net)
$ws.Range("L42").Value = 6506
$ws.Range("M42").Value = -4414
$ws.Range("N42").Value = 202312
$ws.Range("O42").Value = 20338

# Row 43 - External (Gross)
$ws.Range("L43").Value = 8063
$ws.Range("M43").Value = 6989
$ws.Range("N43").Value = 221983
$ws.Range("O43").Value = 61800

# Row 44 - Less: Amortization
$ws.Range("L44").Value = 1557
$ws.Range("M44").Value = 11403
$ws.Range("N44").Value = 19671
$ws.Range("O44").Value = 41462

# Row 47 - Domestic (Gross)
$ws.Range("L47").Value = 180602
$ws.Range("M47").Value = 167045
$ws.Range("N47").Value = 145200
$ws.Range("O47").Value = 67464

# Row 49 - Amortization
$ws.Range("L49").Value = 185
$ws.Range("M49").Value = 122034
$ws.Range("N49").Value = 87
$ws.Range("O49").Value = 120000

# ---------------------------------------------------------------------
# 4) View state: scroll/selection
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 29
$ws.Range("K3").Select()
